$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "model_6_1_1"
$ws.Cells.Item(2,2).Value = 0.05024796887061866
$ws.Cells.Item(2,3).Value = -0.4460531622001842
$ws.Cells.Item(2,4).Value = -3.685831778419148
$ws.Cells.Item(2,5).Value = -0.7497896818492398
$ws.Cells.Item(2,6).Value = 1.051095724105835
$ws.Cells.Item(2,7).Value = 2.425705909729004
$ws.Cells.Item(2,8).Value = 1.069767832756042
$ws.Cells.Item(2,9).Value = 1.787616610527039

$ws.Cells.Item(3,1).Value = "model_6_1_3"
$ws.Cells.Item(3,2).Value = 0.06636021522386093
$ws.Cells.Item(3,3).Value = -0.4947018566597037
$ws.Cells.Item(3,4).Value = -3.487343364435751
$ws.Cells.Item(3,5).Value = -0.7712068475366818
$ws.Cells.Item(3,6).Value = 1.03326416015625
$ws.Cells.Item(3,7).Value = 2.507312774658203
$ws.Cells.Item(3,8).Value = 1.024453282356262
$ws.Cells.Item(3,9).Value = 1.809496760368347

$ws.Cells.Item(4,1).Value = "model_6_1_4"
$ws.Cells.Item(4,2).Value = 0.07368823630507115
$ws.Cells.Item(4,3).Value = -0.5046170790658033
$ws.Cells.Item(4,4).Value = -3.429284598149728
$ws.Cells.Item(4,5).Value = -0.7737205692288598
$ws.Cells.Item(4,6).Value = 1.025154232978821
$ws.Cells.Item(4,7).Value = 2.523944854736328
$ws.Cells.Item(4,8).Value = 1.0111985206604
$ws.Cells.Item(4,9).Value = 1.812064647674561

$ws.Cells.Item(5,1).Value = "model_6_1_2"
$ws.Cells.Item(5,2).Value = 0.0764462544829333
$ws.Cells.Item(5,3).Value = -0.4239671110629386
$ws.Cells.Item(5,4).Value = -3.432375602940859
$ws.Cells.Item(5,5).Value = -0.7039375129460101
$ws.Cells.Item(5,6).Value = 1.022101879119873
$ws.Cells.Item(5,7).Value = 2.388657331466675
$ws.Cells.Item(5,8).Value = 1.011904120445251
$ws.Cells.Item(5,9).Value = 1.74077320098877

$ws.Cells.Item(6,1).Value = "model_6_1_5"
$ws.Cells.Item(6,2).Value = 0.171376367955839
$ws.Cells.Item(6,3).Value = -0.3160953629476637
$ws.Cells.Item(6,4).Value = -2.716821805375643
$ws.Cells.Item(6,5).Value = -0.5349189356743507
$ws.Cells.Item(6,6).Value = 0.9170421957969666
$ws.Cells.Item(6,7).Value = 2.207705974578857
$ws.Cells.Item(6,8).Value = 0.8485444188117981
$ws.Cells.Item(6,9).Value = 1.568100810050964

$ws.Cells.Item(7,1).Value = "model_6_1_17"
$ws.Cells.Item(7,2).Value = 0.1801254719301683
$ws.Cells.Item(7,3).Value = -0.3237610291215822
$ws.Cells.Item(7,4).Value = -2.696788112347436
$ws.Cells.Item(7,5).Value = -0.5394795522406246
$ws.Cells.Item(7,6).Value = 0.9073595404624939
$ws.Cells.Item(7,7).Value = 2.220564842224121
$ws.Cells.Item(7,8).Value = 0.843970775604248
$ws.Cells.Item(7,9).Value = 1.572759866714478

$ws.Cells.Item(8,1).Value = "model_6_1_22"
$ws.Cells.Item(8,2).Value = 0.1805563294075355
$ws.Cells.Item(8,3).Value = -0.3300861536875999
$ws.Cells.Item(8,4).Value = -2.639245461350396
$ws.Cells.Item(8,5).Value = -0.5389230586495828
$ws.Cells.Item(8,6).Value = 0.9068827033042908
$ws.Cells.Item(8,7).Value = 2.231175422668457
$ws.Cells.Item(8,8).Value = 0.8308337926864624
$ws.Cells.Item(8,9).Value = 1.572191596031189

$ws.Cells.Item(9,1).Value = "model_6_1_23"
$ws.Cells.Item(9,2).Value = 0.1809502212297278
$ws.Cells.Item(9,3).Value = -0.3284735625829664
$ws.Cells.Item(9,4).Value = -2.639383797789864
$ws.Cells.Item(9,5).Value = -0.5375364781833325
$ws.Cells.Item(9,6).Value = 0.9064467549324036
$ws.Cells.Item(9,7).Value = 2.228470087051392
$ws.Cells.Item(9,8).Value = 0.8308653831481934
$ws.Cells.Item(9,9).Value = 1.570774912834167

$ws.Cells.Item(10,1).Value = "model_6_1_24"
$ws.Cells.Item(10,2).Value = 0.1811612888872539
$ws.Cells.Item(10,3).Value = -0.3277082246942489
$ws.Cells.Item(10,4).Value = -2.638649588768884
$ws.Cells.Item(10,5).Value = -0.5367966598798979
$ws.Cells.Item(10,6).Value = 0.906213104724884
$ws.Cells.Item(10,7).Value = 2.227186441421509
$ws.Cells.Item(10,8).Value = 0.830697774887085
$ws.Cells.Item(10,9).Value = 1.570019245147705

$ws.Cells.Item(11,1).Value = "model_6_1_21"
$ws.Cells.Item(11,2).Value = 0.1820210889903224
$ws.Cells.Item(11,3).Value = -0.3237083813769415
$ws.Cells.Item(11,4).Value = -2.641832132974327
$ws.Cells.Item(11,5).Value = -0.5336530555740178
$ws.Cells.Item(11,6).Value = 0.9052616357803345
$ws.Cells.Item(11,7).Value = 2.220476627349854
$ws.Cells.Item(11,8).Value = 0.831424355506897
$ws.Cells.Item(11,9).Value = 1.566807627677917

$ws.Cells.Item(12,1).Value = "model_6_1_20"
$ws.Cells.Item(12,2).Value = 0.1826904300514818
$ws.Cells.Item(12,3).Value = -0.3211293366246437
$ws.Cells.Item(12,4).Value = -2.638639188123275
$ws.Cells.Item(12,5).Value = -0.5310748509320187
$ws.Cells.Item(12,6).Value = 0.9045209288597107
$ws.Cells.Item(12,7).Value = 2.216150283813477
$ws.Cells.Item(12,8).Value = 0.8306953907012939
$ws.Cells.Item(12,9).Value = 1.564173579216003

$ws.Cells.Item(13,1).Value = "model_6_1_13"
$ws.Cells.Item(13,2).Value = 0.1829352961531033
$ws.Cells.Item(13,3).Value = -0.3078440340895363
$ws.Cells.Item(13,4).Value = -2.704351071855781
$ws.Cells.Item(13,5).Value = -0.5264346869532492
$ws.Cells.Item(13,6).Value = 0.9042497873306274
$ws.Cells.Item(13,7).Value = 2.193864822387695
$ws.Cells.Item(13,8).Value = 0.8456974029541016
$ws.Cells.Item(13,9).Value = 1.559433102607727

$ws.Cells.Item(14,1).Value = "model_6_1_18"
$ws.Cells.Item(14,2).Value = 0.1830247481214229
$ws.Cells.Item(14,3).Value = -0.3162498167754899
$ws.Cells.Item(14,4).Value = -2.66372317870195
$ws.Cells.Item(14,5).Value = -0.5294730132158547
$ws.Cells.Item(14,6).Value = 0.9041508436203003
$ws.Cells.Item(14,7).Value = 2.207965135574341
$ws.Cells.Item(14,8).Value = 0.8364220857620239
$ws.Cells.Item(14,9).Value = 1.56253719329834

$ws.Cells.Item(15,1).Value = "model_6_1_12"
$ws.Cells.Item(15,2).Value = 0.1831763393842442
$ws.Cells.Item(15,3).Value = -0.3066586574536934
$ws.Cells.Item(15,4).Value = -2.696943288812197
$ws.Cells.Item(15,5).Value = -0.5246267089414367
$ws.Cells.Item(15,6).Value = 0.9039831161499023
$ws.Cells.Item(15,7).Value = 2.191876411437988
$ws.Cells.Item(15,8).Value = 0.8440061807632446
$ws.Cells.Item(15,9).Value = 1.557586073875427

$ws.Cells.Item(16,1).Value = "model_6_1_9"
$ws.Cells.Item(16,2).Value = 0.1836226319111175
$ws.Cells.Item(16,3).Value = -0.2932642844674602
$ws.Cells.Item(16,4).Value = -2.73506402925061
$ws.Cells.Item(16,5).Value = -0.5169894481319772
$ws.Cells.Item(16,6).Value = 0.9034891128540039
$ws.Cells.Item(16,7).Value = 2.169407844543457
$ws.Cells.Item(16,8).Value = 0.8527091145515442
$ws.Cells.Item(16,9).Value = 1.54978358745575

$ws.Cells.Item(17,1).Value = "model_6_1_19"
$ws.Cells.Item(17,2).Value = 0.1837329884844883
$ws.Cells.Item(17,3).Value = -0.3148548572734218
$ws.Cells.Item(17,4).Value = -2.653834964158423
$ws.Cells.Item(17,5).Value = -0.5272134243369704
$ws.Cells.Item(17,6).Value = 0.9033670425415039
$ws.Cells.Item(17,7).Value = 2.205625295639038
$ws.Cells.Item(17,8).Value = 0.8341646194458008
$ws.Cells.Item(17,9).Value = 1.560228586196899

$ws.Cells.Item(18,1).Value = "model_6_1_16"
$ws.Cells.Item(18,2).Value = 0.1839558144942786
$ws.Cells.Item(18,3).Value = -0.3112418940612118
$ws.Cells.Item(18,4).Value = -2.658407115147372
$ws.Cells.Item(18,5).Value = -0.5245565762630227
$ws.Cells.Item(18,6).Value = 0.9031205177307129
$ws.Cells.Item(18,7).Value = 2.199564695358276
$ws.Cells.Item(18,8).Value = 0.8352084159851074
$ws.Cells.Item(18,9).Value = 1.557514429092407

$ws.Cells.Item(19,1).Value = "model_6_1_8"
$ws.Cells.Item(19,2).Value = 0.1845399773884882
$ws.Cells.Item(19,3).Value = -0.2896587808053863
$ws.Cells.Item(19,4).Value = -2.697964365251934
$ws.Cells.Item(19,5).Value = -0.5099551818482313
$ws.Cells.Item(19,6).Value = 0.9024739265441895
$ws.Cells.Item(19,7).Value = 2.163359642028809
$ws.Cells.Item(19,8).Value = 0.8442392945289612
$ws.Cells.Item(19,9).Value = 1.542597413063049

$ws.Cells.Item(20,1).Value = "model_6_1_15"
$ws.Cells.Item(20,2).Value = 0.1846352252023525
$ws.Cells.Item(20,3).Value = -0.3073394395036146
$ws.Cells.Item(20,4).Value = -2.665897531027434
$ws.Cells.Item(20,5).Value = -0.5219510050697309
$ws.Cells.Item(20,6).Value = 0.9023685455322266
$ws.Cells.Item(20,7).Value = 2.193018436431885
$ws.Cells.Item(20,8).Value = 0.8369184732437134
$ws.Cells.Item(20,9).Value = 1.554852485656738

$ws.Cells.Item(21,1).Value = "model_6_1_11"
$ws.Cells.Item(21,2).Value = 0.1853446293840699
$ws.Cells.Item(21,3).Value = -0.2964781794211391
$ws.Cells.Item(21,4).Value = -2.696097202210205
$ws.Cells.Item(21,5).Value = -0.5156852033039667
$ws.Cells.Item(21,6).Value = 0.9015833735466003
$ws.Cells.Item(21,7).Value = 2.174798727035522
$ws.Cells.Item(21,8).Value = 0.8438130021095276
$ws.Cells.Item(21,9).Value = 1.54845130443573

$ws.Cells.Item(22,1).Value = "model_6_1_14"
$ws.Cells.Item(22,2).Value = 0.1853898908969185
$ws.Cells.Item(22,3).Value = -0.3035967241769189
$ws.Cells.Item(22,4).Value = -2.669166528627525
$ws.Cells.Item(22,5).Value = -0.5190416293286011
$ws.Cells.Item(22,6).Value = 0.9015334248542786
$ws.Cells.Item(22,7).Value = 2.186739921569824
$ws.Cells.Item(22,8).Value = 0.837664783000946
$ws.Cells.Item(22,9).Value = 1.551880121231079

$ws.Cells.Item(23,1).Value = "model_6_1_10"
$ws.Cells.Item(23,2).Value = 0.1875309698781547
$ws.Cells.Item(23,3).Value = -0.2859181276612488
$ws.Cells.Item(23,4).Value = -2.691300693772672
$ws.Cells.Item(23,5).Value = -0.5060036862849473
$ws.Cells.Item(23,6).Value = 0.8991637825965881
$ws.Cells.Item(23,7).Value = 2.157084465026855
$ws.Cells.Item(23,8).Value = 0.8427180051803589
$ws.Cells.Item(23,9).Value = 1.538560271263123

$ws.Cells.Item(24,1).Value = "model_6_1_7"
$ws.Cells.Item(24,2).Value = 0.188062533235088
$ws.Cells.Item(24,3).Value = -0.278826496334124
$ws.Cells.Item(24,4).Value = -2.619082946053535
$ws.Cells.Item(24,5).Value = -0.4922435605820743
$ws.Cells.Item(24,6).Value = 0.898575484752655
$ws.Cells.Item(24,7).Value = 2.145188808441162
$ws.Cells.Item(24,8).Value = 0.8262307643890381
$ws.Cells.Item(24,9).Value = 1.524502754211426

$ws.Cells.Item(25,1).Value = "model_6_1_6"
$ws.Cells.Item(25,2).Value = 0.1893593178598063
$ws.Cells.Item(25,3).Value = -0.2712199711673999
$ws.Cells.Item(25,4).Value = -2.584673131951025
$ws.Cells.Item(25,5).Value = -0.4820139109860273
$ws.Cells.Item(25,6).Value = 0.8971403837203979
$ws.Cells.Item(25,7).Value = 2.132429122924805
$ws.Cells.Item(25,8).Value = 0.8183751106262207
$ws.Cells.Item(25,9).Value = 1.514051914215088

$ws.Cells.Item(26,1).Value = "model_6_1_0"
$ws.Cells.Item(26,2).Value = 0.312700901608728
$ws.Cells.Item(26,3).Value = 0.3658468963373929
$ws.Cells.Item(26,4).Value = -0.6676155347333375
$ws.Cells.Item(26,5).Value = 0.273376440866874
$ws.Cells.Item(26,6).Value = 0.7606375813484192
$ws.Cells.Item(26,7).Value = 1.063770771026611
$ws.Cells.Item(26,8).Value = 0.380713939666748
$ws.Cells.Item(26,9).Value = 0.7423316836357117
